# [IMP] re-arrange data mapping to templates
#
# Re-lays out the "Budget Summary Report" sheet: a new title/filter block
# (report title, "ศูนย์"/"จากวันที่"/"ถึงวันที่" filter rows) is inserted
# above the existing column-header row, the header row itself gains a
# "รายการ" column and its contract start/end date columns are renamed and
# swapped, and the sheet-wide look (fonts, wrap text, row heights, column
# widths) is refreshed to match the new template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Blank out the old "filter" rows (old A1 "NSTDA", old A2 report title,
#    old A3/B3 placeholders) so we can rebuild them with the new layout.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).ClearContents()
$ws.Rows.Item(2).ClearContents()
$ws.Rows.Item(3).ClearContents()

# ---------------------------------------------------------------------
# 2. New rows 1-5: report title + filter/criteria block.
# ---------------------------------------------------------------------

# Row 1: report title, merged across A1:B1, bold 11pt.
$ws.Range("A1").Value = "รายงานการจ้างเหมาปฏิบัติงาน"
$ws.Range("A1:B1").Merge()
$ws.Rows.Item(1).RowHeight = 26
$ws.Range("A1:J1").Font.Name = "Arial"
$ws.Range("A1:J1").WrapText = $true
$ws.Range("A1:B1").Font.Size = 11
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("C1:J1").Font.Size = 10

# Row 2: "ศูนย์" (Center) filter line.
$ws.Range("A2").Value = "ศูนย์"
$ws.Rows.Item(2).RowHeight = 13
$ws.Range("A2:J2").Font.Name = "Arial"
$ws.Range("A2:J2").Font.Size = 10
$ws.Range("A2:J2").WrapText = $true
$ws.Range("A2:B2").Font.Bold = $true

# Row 3: "จากวันที่" (From date) filter line.
$ws.Range("A3").Value = "จากวันที่"
$ws.Rows.Item(3).RowHeight = 13
$ws.Range("A3:J3").Font.Name = "Arial"
$ws.Range("A3:J3").Font.Size = 10
$ws.Range("A3:J3").WrapText = $true
$ws.Range("A3").Font.Bold = $true

# Row 4 (new): "ถึงวันที่" (To date) filter line.
$ws.Range("A4").Value = "ถึงวันที่"
$ws.Rows.Item(4).RowHeight = 13
$ws.Range("A4:J4").Font.Name = "Arial"
$ws.Range("A4:J4").Font.Size = 10
$ws.Range("A4:J4").WrapText = $true
$ws.Range("A4").Font.Bold = $true

# Row 5 (new): spacer line before the data table header.
$ws.Rows.Item(5).RowHeight = 13
$ws.Range("A5:J5").Font.Name = "Arial"
$ws.Range("A5:J5").Font.Size = 10
$ws.Range("A5:J5").Font.Bold = $true
$ws.Range("A5:J5").WrapText = $true

# ---------------------------------------------------------------------
# 3. Row 6: data-table header row (re-ordered / re-labelled columns).
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "ลำดับ"
$ws.Range("B6").Value = "ลงวันที่ "
$ws.Range("C6").Value = "เลขที่ สัญญา"
$ws.Range("D6").Value = "เลขที่ PO"
$ws.Range("E6").Value = "ผู้ขาย"
$ws.Range("F6").Value = "รายการ"
$ws.Range("G6").Value = "วันที่เริ่มต้นสัญญา"
$ws.Range("H6").Value = "วันที่สิ้นสุดสัญญา"
$ws.Range("I6").Value = "จำนวนเงิน"
$ws.Range("J6").Value = "สกุลเงิน"

$ws.Rows.Item(6).RowHeight = 16
$ws.Range("A6:J6").Font.Name = "Arial"
$ws.Range("A6:J6").Font.Size = 10
$ws.Range("A6:J6").WrapText = $true
$ws.Range("A6:J6").HorizontalAlignment = -4108
$ws.Range("A6:J6").Borders.LineStyle = 1
$ws.Range("A6:J6").Borders.Weight = 2

# ---------------------------------------------------------------------
# 4. Column widths to roughly match the new template proportions.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 18.8
$ws.Columns.Item(4).ColumnWidth = 18.8
$ws.Columns.Item(5).ColumnWidth = 26.9
$ws.Columns.Item(6).ColumnWidth = 37.1
$ws.Columns.Item(7).ColumnWidth = 16.6
$ws.Columns.Item(8).ColumnWidth = 16.6
$ws.Columns.Item(9).ColumnWidth = 16.6
$ws.Columns.Item(10).ColumnWidth = 16.6

# ---------------------------------------------------------------------
# 5. Misc sheet-view housekeeping.
# ---------------------------------------------------------------------
$ws.Range("E21").Select()
